# Implement vendor and PO approval dialogs with robust timing and data validation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (E, G, Q) ---
# ColumnWidth is expressed in "characters"; the engine applies the standard
# Excel padding (+5/6 of a character) when converting to the persisted XML
# width, so back that padding out to land exactly on the target integers.
$pad = 5 / 6
$ws.Columns.Item(5).ColumnWidth = 13 - $pad    # E: 18 -> 13
$ws.Columns.Item(7).ColumnWidth = 16 - $pad    # G: 19 -> 16
$ws.Columns.Item(17).ColumnWidth = 9 - $pad    # Q: 26 -> 9

# --- Row 2 data update: vendor switched from Johnstone Supply to Slakey ---
# Helper: write a value while forcing the "Text" number format so the engine
# keeps the literal string instead of auto-coercing it to a date/number
# (needed for values that look numeric, e.g. "07/28/25" or "74.20"), then
# reset the cell's style back to Normal so no stray style index lingers.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "B2" "07/28/25"
$ws.Range("D2").Value = "SLABRO"
$ws.Range("E2").Value = "Slakey "
Set-TextValue "G2" "860166732"
Set-TextValue "H2" "07/28/25"
Set-TextValue "I2" "81.53"
Set-TextValue "J2" "7.33"
$ws.Range("K2").ClearContents()
Set-TextValue "L2" "74.20"
Set-TextValue "P2" "1200"
$ws.Range("Q2").Value = "Shop"
$ws.Range("R2").Value = "doc11181820250804113622-5_1754581659155.pdf"
